# Updated cryptos list on Thu Nov  2 09:48:39 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=1 (rank), B=2 (Coin), C=3 (Link), D=4 (Price), E=5 (Volume(1h))

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "35.527.99"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +2.99%  "

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.840.00"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +1.82%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = "  +0.28%  "

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "231.55"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.83%  "

# Row 6 - XRP
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.612"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.19%  "

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = "  +0.22%  "

# Row 8 - Solana
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "43.73"
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +12.27%  "

# Row 9 - Cardano
$ws.Cells.Item(9, 5).Value = "  +7.66%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10, 5).Value = "  +4.98%  "

# Row 11 - TRON
$ws.Cells.Item(11, 5).Value = "  +2.38%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "2.104.98"
$ws.Cells.Item(12, 4).Style = "Normal"

# Row 13 - WrappedEther (was Polygon)
$ws.Cells.Item(13, 2).Value = "WrappedEther"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "1.839.08"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +1.92%  "

# Row 14 - Polygon (was WrappedEther)
$ws.Cells.Item(14, 2).Value = "Polygon"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.675"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +7.04%  "

# Row 15 - Chainlink
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "11.25"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +1.58%  "

# Row 16 - Polkadot
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "4.72"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +7.66%  "

# Row 17 - WrappedBTC
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "35.478.10"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +2.89%  "

# Row 18 - Litecoin
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "70.31"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +3.02%  "

# Row 19 - ShibaInu
$ws.Cells.Item(19, 5).Value = "  +4.16%  "

# Row 20 - BitcoinCash
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "244.43"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.67%  "

# Row 21 - Avalanche
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "12.10"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +8.40%  "

# Row 22 - Uniswap
$ws.Cells.Item(22, 5).Value = "  +14.64%  "

# Row 23 - Dai
$ws.Cells.Item(23, 5).Value = "  +0.28%  "

# Row 24 - Toncoin
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.21"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.20%  "

# Row 25 - Monero
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "172.00"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.55%  "

# Row 26 - Cosmos
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "7.95"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +3.28%  "

# Row 27 - EthereumClassic
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "17.81"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.91%  "

# Row 28 - Stellar
$ws.Cells.Item(28, 5).Value = "  -0.76%  "

# Row 29 - PancakeSwap
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.59"
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +29.64%  "

# Row 30 - BinanceUSD
$ws.Cells.Item(30, 5).Value = "  +0.22%  "

# Row 31 - EURNeutrino
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "3.323.71"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +36.80%  "

# Row 32 - Hedera
$ws.Cells.Item(32, 5).Value = "  +7.38%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Cells.Item(33, 5).Value = "  +6.38%  "

# Row 34 - Filecoin
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "3.95"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.82%  "

# Row 35 - LidoDAOToken
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.86"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +1.63%  "

# Row 36 - Aave
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "96.38"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +17.05%  "

# Row 37 - ImmutableX
$ws.Cells.Item(37, 5).Value = "  +7.48%  "

# Row 38 - TrustWalletToken
$ws.Cells.Item(38, 5).Value = "  +6.05%  "

# Row 39 - Maker
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.351.68"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +3.34%  "

# Row 40 - InjectiveProtocol
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "15.58"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +11.40%  "

# Row 41 - RenderToken
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.46"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +6.26%  "

# Row 42 - VeChain
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0196"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +5.12%  "

# Row 43 - ARBITRUM
$ws.Cells.Item(43, 5).Value = "  +6.27%  "

# Row 44 - WEMIXToken
$ws.Cells.Item(44, 5).Value = "  +4.37%  "

# Row 45 - HuobiToken
$ws.Cells.Item(45, 5).Value = "  +0.92%  "

# Row 46 - MXToken
$ws.Cells.Item(46, 5).Value = "  +0.89%  "

# Row 47 - FraxShare
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "6.28"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +8.48%  "

# Row 48 - Kaspa
$ws.Cells.Item(48, 5).Value = "  +1.16%  "

# Row 49 - RocketPoolETH
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.007.68"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +1.95%  "

# Row 50 - PaxDollar
$ws.Cells.Item(50, 5).Value = "  +0.24%  "

# Row 51 - Quant
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "103.51"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.72%  "
